$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert new "Std ..." columns around the existing ones ---
# A1 Location, B1 Baseline NPV stay put.
$ws.Range("E1").Value = "Opt NPV"
$ws.Range("H1").Value = "Delta NPV"
$ws.Range("C1").Value = "Std baseline NPV"
$ws.Range("F1").Value = "Std opt NPV"
$ws.Range("D1").Value = "Std baseline NPV frac"
$ws.Range("G1").Value = "Std opt NPV frac"
$ws.Range("I1").Value = "Std delta NPV (%)"

# --- Row 2: Braidwood ---
$ws.Range("B2").Value = 2102483898.8800001
$ws.Range("C2").Value = 10442175
$ws.Range("D2").Formula = "=C2/B2"
$ws.Range("E2").Value = 3932705714.9400001
$ws.Range("F2").Value = 46329.74
$ws.Range("G2").Formula = "=F2/E2"
$ws.Range("H2").Formula = "=E2-B2"
$ws.Range("I2").Formula = "=100*SQRT(POWER(G2,2)+POWER(D2,2))"

# --- Row 3: Prairie Island ---
$ws.Range("B3").Value = 262560457.34999999
$ws.Range("C3").Value = 0
$ws.Range("D3").Formula = "=C3/B3"
$ws.Range("E3").Value = 1534342635.98
$ws.Range("F3").Value = 21276.57
$ws.Range("G3").Formula = "=F3/E3"
$ws.Range("H3").Formula = "=E3-B3"
$ws.Range("I3").Formula = "=100*SQRT(POWER(G3,2)+POWER(D3,2))"

# --- Row 4: Davis Besse ---
$ws.Range("B4").Value = 1765764570.5999999
$ws.Range("C4").Value = 0

# --- Row 5: Houston ---
$ws.Range("B5").Value = 2861165724.4899998
$ws.Range("C5").Value = 75472956.379999995

# --- Row 6: Cooper ---
$ws.Range("B6").Value = 1080996406.4100001
$ws.Range("C6").Value = 0

# --- Column widths (closest reachable to the recorded 12.1640625 / 11.1640625) ---
$ws.Columns.Item(3).ColumnWidth = 11.33
$ws.Columns.Item(4).ColumnWidth = 11.33
$ws.Columns.Item(5).ColumnWidth = 10.33
$ws.Columns.Item(7).ColumnWidth = 11.33
$ws.Columns.Item(8).ColumnWidth = 10.33

# --- Selection moves to J5 ---
$ws.Range("J5").Select() | Out-Null
